$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: insert 8 new data rows (20:27) using row 19's formatting
#    (a "normal" interior data row, not the bottom-border closing row).
#    This pushes the old row 20 (closing-border row, JHON FREDY's record)
#    down to row 28, and the footer rows (25/26) down to rows 33/34 --
#    matching the target dimension B2:J34 and merged cells.
# ---------------------------------------------------------------------------
$ws.Rows("20:27").Insert()
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J27").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Update the summary block
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 502320     # VALOR MORA (total)
$ws.Range("C13").Value = 9          # Cant. Trabajadores
$ws.Range("F13").Value = 3          # Cant. Periodos

# ---------------------------------------------------------------------------
# 3) Rewrite the full worker/period data block (rows 16-28)
#    Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#             E=Periodo Mora, F=Valor Mora, G=Salario Basico
# ---------------------------------------------------------------------------
$data = @(
    @("CC", "94483278",   "EIDER PIMENTEL CALDON",           "1606", 34000, 0),
    @("CC", "94483278",   "EIDER PIMENTEL CALDON",           "1605", 34000, 0),
    @("CC", "1047427860", "YINETH PAOLA BALLESTAS FERIA",    "1606", 32000, 0),
    @("CC", "1047427860", "YINETH PAOLA BALLESTAS FERIA",    "1605", 32000, 0),
    @("CC", "1116435458", "JORGE ALBERTO SANTAMARIA",        "1612", 35600, 890000),
    @("CC", "1118288813", "MARIA FERNANDA ORTEGA OSPINA",    "1612", 62720, 1568000),
    @("CC", "22790491",   "EMILDA ISABEL ATENCIA PALMERA",   "1612", 26000, 650000),
    @("CC", "14899808",   "JHON FREDY GUTIERREZ CASTILLO",   "1612", 68000, 1700000),
    @("CC", "1047446376", "ANGEL DE JESUS MESTRA ZULETA",    "1612", 34000, 850000),
    @("CC", "1143338624", "PEDRO LUIS PULIDO ANAYA",         "1606", 36000, 0),
    @("CC", "1143338624", "PEDRO LUIS PULIDO ANAYA",         "1605", 36000, 0),
    @("CC", "80812126",   "HENIEK YITZAK HERRERA RAMIREZ",   "1606", 36000, 0),
    @("CC", "80812126",   "HENIEK YITZAK HERRERA RAMIREZ",   "1605", 36000, 0)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $row[1]   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[2]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[3]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[4]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[5]   # G - Salario Basico
}

Write-Output "Edit complete"
